$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.826.24'
$ws.Range("E2").Value = '  -0.32%  '
$ws.Range("D3").Value = '3.505.04'
$ws.Range("E3").Value = '  -0.78%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.86'
$ws.Range("E5").Value = '  -0.38%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.56'
$ws.Range("E6").Value = '  -2.42%  '
$ws.Range("D7").Value = '3.504.31'
$ws.Range("E7").Value = '  -0.75%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.494'
$ws.Range("E9").Value = '  +0.35%  '
$ws.Range("E10").Value = '  -0.39%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.11'
$ws.Range("E11").Value = '  +3.24%  '
$ws.Range("E12").Value = '  -0.68%  '
$ws.Range("D13").Value = '4.108.97'
$ws.Range("E13").Value = '  -0.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.06'
$ws.Range("E14").Value = '  -0.65%  '
$ws.Range("E15").Value = '  -0.74%  '
$ws.Range("B16").Value = 'TRON'
$ws.Range("C16").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.117'
$ws.Range("E16").Value = '  -0.25%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.507.61'
$ws.Range("E17").Value = '  -0.82%  '
$ws.Range("D18").Value = '64.869.26'
$ws.Range("E18").Value = '  -0.29%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.72'
$ws.Range("E19").Value = '  -3.87%  '
$ws.Range("E20").Value = '  +1.45%  '
$ws.Range("E21").Value = '  -2.59%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '388.34'
$ws.Range("E22").Value = '  -0.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.576'
$ws.Range("E23").Value = '  +0.64%  '
$ws.Range("B24").Value = 'WrappedeETH'
$ws.Range("C24").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D24").Value = '3.648.98'
$ws.Range("E24").Value = '  -0.83%  '
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.18'
$ws.Range("E25").Value = '  +0.47%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000113'
$ws.Range("E27").Value = '  +0.67%  '
$ws.Range("B28").Value = 'Fetch.AI'
$ws.Range("C28").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.61'
$ws.Range("E28").Value = '  +17.56%  '
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.72'
$ws.Range("E29").Value = '  -0.48%  '
$ws.Range("E30").Value = '  +0.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.27'
$ws.Range("E31").Value = '  +1.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.36'
$ws.Range("E32").Value = '  +1.75%  '
$ws.Range("D33").Value = '3.514.47'
$ws.Range("E33").Value = '  -1.07%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '24.02'
$ws.Range("E34").Value = '  +0.95%  '
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.145'
$ws.Range("E36").Value = '  +0.81%  '
$ws.Range("E37").Value = '  +5.27%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '169.77'
$ws.Range("E38").Value = '  +0.50%  '
$ws.Range("E39").Value = '  -0.91%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.81'
$ws.Range("E40").Value = '  -1.21%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0825'
$ws.Range("E41").Value = '  +3.22%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.815'
$ws.Range("E42").Value = '  -0.73%  '
$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '42.52'
$ws.Range("E43").Value = '  +0.27%  '
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  -0.07%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '25.41'
$ws.Range("E45").Value = '  -3.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.22'
$ws.Range("E46").Value = '  +2.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.41'
$ws.Range("E47").Value = '  -0.08%  '
$ws.Range("E48").Value = '  -1.60%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.91'
$ws.Range("E49").Value = '  +1.27%  '
$ws.Range("D50").Value = '2.359.80'
$ws.Range("E50").Value = '  -1.14%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0267'
$ws.Range("E51").Value = '  +2.43%  '
